$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1080.5
$ws.Range("J17").Value = 1106
$ws.Range("L17").Value = 3318
$ws.Range("N17").Value = -3654
$ws.Range("H32").Value = 15994.167
$ws.Range("I32").Value = 14989.333
$ws.Range("J32").Value = 16999
$ws.Range("K32").Value = 14989.333
$ws.Range("L32").Value = 16999
$ws.Range("M32").Value = -14663.333
$ws.Range("N32").Value = -17651
$ws.Range("H53").Value = 620.8929000000001
$ws.Range("I53").Value = 739.05554
$ws.Range("K53").Value = 739.05554
$ws.Range("M53").Value = -102.05554
$ws.Range("H57").Value = 68846.5
$ws.Range("J57").Value = 68846.5
$ws.Range("L57").Value = 206539.5
$ws.Range("N57").Value = -207537.5
$ws.Range("H70").Value = 556.5294
$ws.Range("J70").Value = 498
$ws.Range("L70").Value = 1494
$ws.Range("N70").Value = -2034
$ws.Range("H73").Value = 556.5294
$ws.Range("J73").Value = 498
$ws.Range("L73").Value = 1494
$ws.Range("N73").Value = -3366
$ws.Range("H80").Value = 3585.9285
$ws.Range("I80").Value = 611.5
$ws.Range("K80").Value = 1834.5
$ws.Range("M80").Value = -836.5
$ws.Range("H83").Value = 3585.9285
$ws.Range("I83").Value = 611.5
$ws.Range("K83").Value = 5503.5
$ws.Range("M83").Value = -511.5
$ws.Range("H101").Value = 1792.6666
$ws.Range("I101").Value = 689.5
$ws.Range("J101").Value = 3999
$ws.Range("K101").Value = 2068.5
$ws.Range("L101").Value = 11997
$ws.Range("M101").Value = -446.5
$ws.Range("N101").Value = -15241
$ws.Range("H111").Value = 3938.889
$ws.Range("I111").Value = 3795.2
$ws.Range("J111").Value = 4118.5
$ws.Range("K111").Value = 11385.6
$ws.Range("L111").Value = 12355.5
$ws.Range("M111").Value = -8318.599999999999
$ws.Range("N111").Value = -18489.5
$ws.Range("H116").Value = 207596.8
$ws.Range("J116").Value = 504999.5
$ws.Range("L116").Value = 504999.5
$ws.Range("N116").Value = -511883.5
$ws.Range("H132").Value = 2175.0715
$ws.Range("I132").Value = 2104.3333
$ws.Range("J132").Value = 2599.5
$ws.Range("K132").Value = 6312.999899999999
$ws.Range("L132").Value = 7798.5
$ws.Range("M132").Value = -3782.999899999999
$ws.Range("N132").Value = -12858.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 20006
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 20006
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 20006
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = -20294
$ws.Range("H32").Value = 1647.3513
$ws.Range("I32").Value = 1647.3513
$ws.Range("K32").Value = 1647.3513
$ws.Range("M32").Value = -1360.3513
$ws.Range("H88").Value = 4727.273
$ws.Range("I88").Value = 3850
$ws.Range("K88").Value = 3850
$ws.Range("M88").Value = -3444
$ws.Range("H91").Value = 4727.273
$ws.Range("I91").Value = 3850
$ws.Range("K91").Value = 3850
$ws.Range("M91").Value = -2446
$ws.Range("H130").Value = 59998.332
$ws.Range("J130").Value = 59998.332
$ws.Range("L130").Value = 59998.332
$ws.Range("N130").Value = -70038.33199999999
$ws.Range("H138").Value = 75428.5
$ws.Range("J138").Value = 75428.5
$ws.Range("L138").Value = 75428.5
$ws.Range("N138").Value = -85708.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 3942.8
$ws.Range("I11").Value = 1052.25
$ws.Range("K11").Value = 1052.25
$ws.Range("M11").Value = -912.25
$ws.Range("H86").Value = 2224.6843
$ws.Range("I86").Value = 1887.091
$ws.Range("K86").Value = 1887.091
$ws.Range("M86").Value = -764.0909999999999
$ws.Range("H89").Value = 2224.6843
$ws.Range("I89").Value = 1887.091
$ws.Range("K89").Value = 9435.455
$ws.Range("M89").Value = -3819.455
$ws.Range("H99").Value = 3683.5557
$ws.Range("I99").Value = 4025.375
$ws.Range("K99").Value = 4025.375
$ws.Range("M99").Value = -2527.375
$ws.Range("H120").Value = 48000
$ws.Range("J120").Value = 48000
$ws.Range("L120").Value = 48000
$ws.Range("N120").Value = -57676
$ws.Range("H122").Value = 86925
$ws.Range("J122").Value = 86925
$ws.Range("L122").Value = 86925
$ws.Range("N122").Value = -96725

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 223.75
$ws.Range("I7").Value = 198.11111
$ws.Range("J7").Value = 269.9
$ws.Range("K7").Value = 198.11111
$ws.Range("L7").Value = 269.9
$ws.Range("M7").Value = -85.11111
$ws.Range("N7").Value = -495.9
$ws.Range("H18").Value = 28386.25
$ws.Range("J18").Value = 28386.25
$ws.Range("L18").Value = 28386.25
$ws.Range("N18").Value = -28846.25
$ws.Range("H31").Value = 2239.9412
$ws.Range("I31").Value = 1003
$ws.Range("K31").Value = 1003
$ws.Range("M31").Value = -708
$ws.Range("H34").Value = 2239.9412
$ws.Range("I34").Value = 1003
$ws.Range("K34").Value = 1003
$ws.Range("M34").Value = -801
$ws.Range("H122").Value = 2432.3333
$ws.Range("I122").Value = 2378.2
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 7134.599999999999
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4684.599999999999
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 134635.86
$ws.Range("I4").Value = 627.2
$ws.Range("J4").Value = 201640.2
$ws.Range("K4").Value = 1881.6
$ws.Range("L4").Value = 604920.6000000001
$ws.Range("M4").Value = -1769.6
$ws.Range("N4").Value = -605144.6000000001
$ws.Range("H75").Value = 5185.2
$ws.Range("I75").Value = 4308.6665
$ws.Range("K75").Value = 12925.9995
$ws.Range("M75").Value = -11927.9995
$ws.Range("H78").Value = 5185.2
$ws.Range("I78").Value = 4308.6665
$ws.Range("K78").Value = 38777.9985
$ws.Range("M78").Value = -33785.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 20004
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 20004
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 20004
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = -20228
$ws.Range("H132").Value = 1298
$ws.Range("I132").Value = 1006
$ws.Range("J132").Value = 1492.6666
$ws.Range("K132").Value = 3018
$ws.Range("L132").Value = 4477.9998
$ws.Range("M132").Value = -488
$ws.Range("N132").Value = -9537.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1545.75
$ws.Range("I16").Value = 1545.75
$ws.Range("K16").Value = 1545.75
$ws.Range("M16").Value = -1375.75
$ws.Range("H46").Value = 2251.9048
$ws.Range("I46").Value = 1457.4
$ws.Range("J46").Value = 2500.1875
$ws.Range("K46").Value = 1457.4
$ws.Range("L46").Value = 2500.1875
$ws.Range("M46").Value = -1269.4
$ws.Range("N46").Value = -2876.1875
$ws.Range("H136").Value = 45457610
$ws.Range("I136").Value = 2340.4285
$ws.Range("J136").Value = 125004330
$ws.Range("K136").Value = 7021.2855
$ws.Range("L136").Value = 375012990
$ws.Range("M136").Value = -4471.2855
$ws.Range("N136").Value = -375018090
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 23224
$ws.Range("I28").Value = 2105
$ws.Range("J28").Value = 26743.834
$ws.Range("K28").Value = 2105
$ws.Range("L28").Value = 26743.834
$ws.Range("M28").Value = -1757
$ws.Range("N28").Value = -27439.834
$ws.Range("H81").Value = 17014.9
$ws.Range("I81").Value = 16544.625
$ws.Range("J81").Value = 18896
$ws.Range("K81").Value = 33089.25
$ws.Range("L81").Value = 37792
$ws.Range("M81").Value = -32028.25
$ws.Range("N81").Value = -39914
$ws.Range("H84").Value = 17014.9
$ws.Range("I84").Value = 16544.625
$ws.Range("J84").Value = 18896
$ws.Range("K84").Value = 165446.25
$ws.Range("L84").Value = 188960
$ws.Range("M84").Value = -160142.25
$ws.Range("N84").Value = -199568
$ws.Range("H121").Value = 59998.5
$ws.Range("J121").Value = 59998.5
$ws.Range("L121").Value = 59998.5
$ws.Range("N121").Value = -63492.5
$ws.Range("H137").Value = 45000
$ws.Range("J137").Value = 45000
$ws.Range("L137").Value = 45000
